$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1814.5333
$ws.Range("J17").Value = 1814.5333
$ws.Range("L17").Value = 5443.5999
$ws.Range("N17").Value = -5779.5999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 3999.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 3999.5
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -4137.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7325.857
$ws.Range("I62").Value = 2425
$ws.Range("J62").Value = 11001.5
$ws.Range("K62").Value = 2425
$ws.Range("L62").Value = 11001.5
$ws.Range("M62").Value = -1801
$ws.Range("N62").Value = -12249.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 7325.857
$ws.Range("I65").Value = 2425
$ws.Range("J65").Value = 11001.5
$ws.Range("K65").Value = 12125
$ws.Range("L65").Value = 55007.5
$ws.Range("M65").Value = -9005
$ws.Range("N65").Value = -61247.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3689.6667
$ws.Range("I74").Value = 4534.5
$ws.Range("K74").Value = 4534.5
$ws.Range("M74").Value = -3598.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3689.6667
$ws.Range("I77").Value = 4534.5
$ws.Range("K77").Value = 22672.5
$ws.Range("M77").Value = -17992.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 4027.5715
$ws.Range("I106").Value = 4027.5715
$ws.Range("K106").Value = 4027.5715
$ws.Range("M106").Value = -3396.5715

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 6501.25
$ws.Range("I113").Value = 5335
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 5335
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -2081
$ws.Range("N113").Value = -16508

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2436.0454
$ws.Range("I137").Value = 1449.125
$ws.Range("K137").Value = 4347.375
$ws.Range("M137").Value = -1797.375

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2239
$ws.Range("I138").Value = 1270
$ws.Range("J138").Value = 4500
$ws.Range("K138").Value = 3810
$ws.Range("L138").Value = 13500
$ws.Range("M138").Value = 1330
$ws.Range("N138").Value = -23780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4638.8945
$ws.Range("I61").Value = 3475.9333
$ws.Range("K61").Value = 3475.9333
$ws.Range("M61").Value = -3263.9333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5536.25
$ws.Range("I74").Value = 5498.615
$ws.Range("K74").Value = 5498.615
$ws.Range("M74").Value = -4624.615

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 5536.25
$ws.Range("I77").Value = 5498.615
$ws.Range("K77").Value = 27493.075
$ws.Range("M77").Value = -23125.075

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3335.2104
$ws.Range("I102").Value = 1669.2858
$ws.Range("K102").Value = 1669.2858
$ws.Range("M102").Value = -47.28580000000011

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4638.8945
$ws.Range("I136").Value = 3475.9333
$ws.Range("K136").Value = 10427.7999
$ws.Range("M136").Value = -7877.7999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2940.25
$ws.Range("I20").Value = 2391.1667
$ws.Range("K20").Value = 2391.1667
$ws.Range("M20").Value = -2144.1667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4141.3687
$ws.Range("I86").Value = 1762.3636
$ws.Range("K86").Value = 1762.3636
$ws.Range("M86").Value = -639.3635999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4141.3687
$ws.Range("I89").Value = 1762.3636
$ws.Range("K89").Value = 8811.817999999999
$ws.Range("M89").Value = -3195.817999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1879.75
$ws.Range("I16").Value = 1840.3334
$ws.Range("J16").Value = 1998
$ws.Range("K16").Value = 1840.3334
$ws.Range("L16").Value = 1998
$ws.Range("M16").Value = -1553.3334
$ws.Range("N16").Value = -2572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4663.825
$ws.Range("I31").Value = 2151.2593
$ws.Range("K31").Value = 2151.2593
$ws.Range("M31").Value = -1856.2593

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4663.825
$ws.Range("I34").Value = 2151.2593
$ws.Range("K34").Value = 2151.2593
$ws.Range("M34").Value = -1949.2593

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 70000
$ws.Range("J109").Value = 70000
$ws.Range("L109").Value = 70000
$ws.Range("N109").Value = -72080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1879.75
$ws.Range("I113").Value = 1840.3334
$ws.Range("J113").Value = 1998
$ws.Range("K113").Value = 1840.3334
$ws.Range("L113").Value = 1998
$ws.Range("M113").Value = 329.6666
$ws.Range("N113").Value = -6338

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2320.2693
$ws.Range("I134").Value = 2346.375
$ws.Range("K134").Value = 7039.125
$ws.Range("M134").Value = -4504.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 703333.7
$ws.Range("I4").Value = 1000000.5
$ws.Range("K4").Value = 3000001.5
$ws.Range("M4").Value = -2999889.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2329.1428
$ws.Range("I80").Value = 2174.5
$ws.Range("J80").Value = 2535.3333
$ws.Range("K80").Value = 2174.5
$ws.Range("L80").Value = 2535.3333
$ws.Range("M80").Value = -1176.5
$ws.Range("N80").Value = -4531.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2329.1428
$ws.Range("I83").Value = 2174.5
$ws.Range("J83").Value = 2535.3333
$ws.Range("K83").Value = 10872.5
$ws.Range("L83").Value = 12676.6665
$ws.Range("M83").Value = -5880.5
$ws.Range("N83").Value = -22660.6665

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4249.5
$ws.Range("I132").Value = 3999.3333
$ws.Range("K132").Value = 11997.9999
$ws.Range("M132").Value = -9467.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6178.143
$ws.Range("I61").Value = 4749.6665
$ws.Range("K61").Value = 4749.6665
$ws.Range("M61").Value = -4547.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 6178.143
$ws.Range("I113").Value = 4749.6665
$ws.Range("K113").Value = 4749.6665
$ws.Range("M113").Value = -2579.6665

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 38332.332
$ws.Range("I27").Value = 34999
$ws.Range("K27").Value = 34999
$ws.Range("M27").Value = -34930

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 810.0714
$ws.Range("I100").Value = 542.44446
$ws.Range("J100").Value = 1291.8
$ws.Range("K100").Value = 1084.88892
$ws.Range("L100").Value = 2583.6
$ws.Range("M100").Value = -543.8889200000001
$ws.Range("N100").Value = -3665.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H115").Value = 29998
$ws.Range("I115").Value = 29998
$ws.Range("K115").Value = 29998
$ws.Range("M115").Value = -28431

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3352.08
$ws.Range("I136").Value = 2065.9333
$ws.Range("K136").Value = 6197.7999
$ws.Range("M136").Value = -3647.7999
